{"js": "// Update the date line and all 100 math-answer cells in the single table.\n// Both were \"regenerated\" in the source data (new random problems/answers);\n// we reproduce the exact target text for every cell, in document order,\n// while leaving all run/paragraph formatting untouched.\n\n// 1) Update the date paragraph text (keeps its run formatting).\nconst dateResults = context.document.body.search(\"2025-04-03 Thursday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2025-04-04 Friday\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Update every cell of the 20x5 answer table (row-major order), keeping\n// each cell's paragraph/run formatting intact (Table.values only touches text).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"29+44=73\", \"82-45=37\", \"34-28=6\", \"49+29=78\", \"27+49=76\"],\n  [\"20-14=6\", \"90-73=17\", \"51-42=9\", \"60-42=18\", \"37+39=76\"],\n  [\"32+19=51\", \"31-18=13\", \"38-9=29\", \"41-22=19\", \"61-54=7\"],\n  [\"90-43=47\", \"21-6=15\", \"24+38=62\", \"73-58=15\", \"48-39=9\"],\n  [\"11-3=8\", \"45-19=26\", \"37+9=46\", \"26+69=95\", \"93-76=17\"],\n  [\"61-37=24\", \"20-11=9\", \"90-23=67\", \"29+5=34\", \"27+48=75\"],\n  [\"92-19=73\", \"27+57=84\", \"90-34=56\", \"87+4=91\", \"36+16=52\"],\n  [\"52+29=81\", \"29+19=48\", \"48+13=61\", \"51-9=42\", \"95-49=46\"],\n  [\"94-19=75\", \"28+9=37\", \"30-9=21\", \"34+17=51\", \"19+32=51\"],\n  [\"70-45=25\", \"93-59=34\", \"18+67=85\", \"55-28=27\", \"23+59=82\"],\n  [\"30-3=27\", \"10-4=6\", \"48+43=91\", \"14+17=31\", \"80-13=67\"],\n  [\"46+7=53\", \"46+7=53\", \"58+25=83\", \"22-6=16\", \"24+27=51\"],\n  [\"19+79=98\", \"53-25=28\", \"9+46=55\", \"7+35=42\", \"49+29=78\"],\n  [\"70-21=49\", \"56+28=84\", \"58+13=71\", \"59+5=64\", \"22-6=16\"],\n  [\"56-47=9\", \"67+4=71\", \"92-59=33\", \"49+4=53\", \"17+24=41\"],\n  [\"95-87=8\", \"25-6=19\", \"20-9=11\", \"93-4=89\", \"27+48=75\"],\n  [\"88-59=29\", \"33-9=24\", \"9+52=61\", \"13+38=51\", \"26+39=65\"],\n  [\"49+5=54\", \"67-59=8\", \"56+19=75\", \"59+6=65\", \"29+12=41\"],\n  [\"26+48=74\", \"56+19=75\", \"29+58=87\", \"53-24=29\", \"72-14=58\"],\n  [\"63-18=45\", \"7+37=44\", \"49+43=92\", \"49+39=88\", \"14+7=21\"],\n];\nawait context.sync();\n", "ps1": "# Update the date line and all 100 math-answer cells in the single table.\n# Both were \"regenerated\" in the source data (new random problems/answers);\n# we reproduce the exact target text for every cell, in document order,\n# while leaving all run/paragraph formatting untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph text (keeps its run formatting).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2025-04-03 Thursday\"\n$find.Replacement.Text = \"2025-04-04 Friday\"\n[void]$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Update every cell of the 20x5 answer table (row-major order), keeping\n# each cell's paragraph/run formatting intact (only Range.Text is replaced).\n$values = @(\n    @(\"29+44=73\", \"82-45=37\", \"34-28=6\", \"49+29=78\", \"27+49=76\"),\n    @(\"20-14=6\", \"90-73=17\", \"51-42=9\", \"60-42=18\", \"37+39=76\"),\n    @(\"32+19=51\", \"31-18=13\", \"38-9=29\", \"41-22=19\", \"61-54=7\"),\n    @(\"90-43=47\", \"21-6=15\", \"24+38=62\", \"73-58=15\", \"48-39=9\"),\n    @(\"11-3=8\", \"45-19=26\", \"37+9=46\", \"26+69=95\", \"93-76=17\"),\n    @(\"61-37=24\", \"20-11=9\", \"90-23=67\", \"29+5=34\", \"27+48=75\"),\n    @(\"92-19=73\", \"27+57=84\", \"90-34=56\", \"87+4=91\", \"36+16=52\"),\n    @(\"52+29=81\", \"29+19=48\", \"48+13=61\", \"51-9=42\", \"95-49=46\"),\n    @(\"94-19=75\", \"28+9=37\", \"30-9=21\", \"34+17=51\", \"19+32=51\"),\n    @(\"70-45=25\", \"93-59=34\", \"18+67=85\", \"55-28=27\", \"23+59=82\"),\n    @(\"30-3=27\", \"10-4=6\", \"48+43=91\", \"14+17=31\", \"80-13=67\"),\n    @(\"46+7=53\", \"46+7=53\", \"58+25=83\", \"22-6=16\", \"24+27=51\"),\n    @(\"19+79=98\", \"53-25=28\", \"9+46=55\", \"7+35=42\", \"49+29=78\"),\n    @(\"70-21=49\", \"56+28=84\", \"58+13=71\", \"59+5=64\", \"22-6=16\"),\n    @(\"56-47=9\", \"67+4=71\", \"92-59=33\", \"49+4=53\", \"17+24=41\"),\n    @(\"95-87=8\", \"25-6=19\", \"20-9=11\", \"93-4=89\", \"27+48=75\"),\n    @(\"88-59=29\", \"33-9=24\", \"9+52=61\", \"13+38=51\", \"26+39=65\"),\n    @(\"49+5=54\", \"67-59=8\", \"56+19=75\", \"59+6=65\", \"29+12=41\"),\n    @(\"26+48=74\", \"56+19=75\", \"29+58=87\", \"53-24=29\", \"72-14=58\"),\n    @(\"63-18=45\", \"7+37=44\", \"49+43=92\", \"49+39=88\", \"14+7=21\")\n)\n\n$tbl = $d.Tables.Item(1)\nfor ($r = 1; $r -le $values.Count; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        [void]($tbl.Cell($r, $c).Range.Text = $row[$c - 1])\n    }\n}\n"}
